$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 342, pushing the existing rows 342-365 down to 346-369.
$ws.Rows.Item(342).Resize(4).Insert()

# Common (constant) metadata shared by every data row in this sheet.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$codreg    = 16
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103006
$categoria   = "Nectarín"
$unidad = "$/caja 16 kilos empedrada"
$origen = "Región de O'Higgins"
$kgUnidad = 16

# New weekly rows (342-345): fecha, variedad, calidad, volumen, precio min, precio max, precio promedio, precio $/Kg
$newRows = @(
    @{ Row = 342; Fecha = 44918; Variedad = "Artic Star";  Calidad = "Primera"; Volumen = 160; Min = 13000; Max = 14000; Prom = 13500; PrecioKg = 844 },
    @{ Row = 343; Fecha = 44918; Variedad = "Artic Star";  Calidad = "Segunda"; Volumen = 80;  Min = 12000; Max = 12000; Prom = 12000; PrecioKg = 750 },
    @{ Row = 344; Fecha = 44918; Variedad = "Super Queen"; Calidad = "Primera"; Volumen = 160; Min = 12000; Max = 13000; Prom = 12500; PrecioKg = 781 },
    @{ Row = 345; Fecha = 44918; Variedad = "Super Queen"; Calidad = "Segunda"; Volumen = 80;  Min = 11000; Max = 11000; Prom = 11000; PrecioKg = 688 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value  = $mercadoId
    $ws.Cells.Item($row, 2).Value  = $mercado
    $ws.Cells.Item($row, 3).Value  = $region
    $ws.Cells.Item($row, 4).Value  = $r.Fecha
    $ws.Cells.Item($row, 5).Value  = $codreg
    $ws.Cells.Item($row, 6).Value  = $tipo
    $ws.Cells.Item($row, 7).Value  = $productoId
    $ws.Cells.Item($row, 8).Value  = $producto
    $ws.Cells.Item($row, 9).Value  = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $r.Variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Min
    $ws.Cells.Item($row, 15).Value = $r.Max
    $ws.Cells.Item($row, 16).Value = $r.Prom
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}
